# Update faturamento_diario_lojas.xlsx with new AA/AB values and recalculated AG totals
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Row 2 - Bibi Cell Mundi
$ws.Range("AA2").Value = 7424.21
$ws.Range("AB2").Value = 1129.2
$ws.Range("AG2").Value = 273326.15

# Row 3 - Bibi Cell Vieiralves
$ws.Range("AA3").Value = 21228.5
$ws.Range("AB3").Value = 55
$ws.Range("AG3").Value = 141813.59

# Row 4 - Bibi Cell Manauara
$ws.Range("AA4").Value = 2569
$ws.Range("AG4").Value = 79963.35000000001

# Row 5 - Bibi Cell Ponta Negra
$ws.Range("AA5").Value = 6123.4
$ws.Range("AB5").Value = 130
$ws.Range("AG5").Value = 74625.00999999999

# Row 6 - total row
$ws.Range("AA6").Value = 37345.11
$ws.Range("AB6").Value = 1314.2
$ws.Range("AG6").Value = 569728.1
